$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "05/01/2026 07:23:00"
$ws.Range("B46").Value = ""
$ws.Range("C46").Value = "VEJA"
$ws.Range("D46").Value = "O aceno a Motta no primeiro projeto enviado por Lula ao Congresso em 2026"
$ws.Range("E46").Value = "https://veja.abril.com.br/coluna/radar/o-primeiro-projeto-de-lei-apresentado-por-lula-ao-congresso-em-2026/"
$ws.Range("F46").Value = "câmara"
$ws.Range("G46").Value = "Proposta enviada pelo governo federal à Câmara cria nova unidade na Rede Federal de Educação Profissional, Científica e Tecnológica"
